$d = $word.ActiveDocument

$targets = @(
    "default values in function definition",
    "2 useful (modern) call-by-references",
    "useful string class usage",
    "useful usage of (modern) file-I/O",
    "useful exception handling"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    foreach ($t in $targets) {
        if ($text -eq $t) {
            $p.Range.HighlightColorIndex = 7
        }
    }
}
